$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.157.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.59%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.031.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.36%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.46%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.26"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.35%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.643"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.45%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.49"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +15.57%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.27%  "

# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.373"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.61%  "

# Row 10
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.95"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.23%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0748"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.68%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.90%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.900"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.38%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.00"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +7.38%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.322.52"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.51"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.74%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.26"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +19.64%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.022.52"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.55%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.916.56"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.09%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.87"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.32%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0870"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.55%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.30"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.51%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.79"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.46%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.78"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +24.83%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.14%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.49"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.43%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.31"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.78%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.77"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.13%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.121"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.52%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.16"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.11%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.21"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.49%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.112"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +26.31%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.69"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +9.88%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0612"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.17%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.44"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +13.39%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.32%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.81%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.98"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +22.57%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.105"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +18.23%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.23"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.43%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.69%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +23.28%  "

# Row 44
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.14"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.68%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0217"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.71%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.03"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +9.91%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.79"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +10.03%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.89"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.30%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.424.92"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.55%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.94"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.04%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.26"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.79%  "
